# The sheet already holds 7 match rows (rows 2-8) for Eoin Morgan.
# The update appends the very same 7 rows again (rows 9-15), this time
# listed in chronological order, extending the table from A1:K8 to A1:K15.
#
# Copying each existing row to its new destination (rather than retyping
# the values) guarantees the appended cells end up with exactly the same
# text/number-as-text formatting and content as their source row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row, matching the chronological order added
# to the bottom of the sheet
$mapping = [ordered]@{
    9  = 7   #  Abu Dhabi    - September 23 2020 - Mumbai Indians
    10 = 3   #  Abu Dhabi    - September 26 2020 - Sunrisers Hyderabad
    11 = 8   #  Dubai (DSC)  - September 30 2020 - Rajasthan Royals
    12 = 2   #  Abu Dhabi    - October 10 2020   - Kings XI Punjab
    13 = 5   #  Sharjah      - October 03 2020   - Delhi Capitals
    14 = 4   #  Sharjah      - October 12 2020   - Royal Challengers Bangalore
    15 = 6   #  Abu Dhabi    - October 07 2020   - Chennai Super Kings
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $ws.Range("A" + $srcRow + ":K" + $srcRow).Copy($ws.Range("A" + $destRow))
}

Write-Output "Appended rows 9-15; used range is now $($ws.UsedRange.Address())"
